# Bug fix in Eduati data files:
# Sheet1 ("HCT116_noCTRL_meas") had a block of leftover/garbage rows
# (45:87, column A only) that do not belong to the data table (which
# really only spans rows 1:44, same as the other two sheets). Remove
# them and refresh the window/selection state to match the fixed file.

$wb = $excel.ActiveWorkbook

# --- Sheet1: drop the stray rows 45:87 (only col A had leftover numbers) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A45:A87").EntireRow.Delete()

# --- Fix up view/selection state ---
# Make Sheet1 the active/selected sheet (was Sheet3 before the fix) and
# move the selection to H40, matching the saved cursor position after
# trimming the bogus rows.
$ws1.Activate()
$ws1.Range("H40").Select()

# Sheet3 is no longer the active sheet; its own selection (A2:N44) is kept
# as-is, only the tabSelected flag moves off of it.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A2:N44").Select()

# Re-activate Sheet1 so it is the tab that is selected/active on save.
$ws1.Activate()
